$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161740779876709
$ws.Range("B1").Value = 2.125133037567139
$ws.Range("C1").Value = 3.39498496055603
$ws.Range("D1").Value = 3.630847454071045
$ws.Range("E1").Value = 1.18119490146637
